# Apply updated NATMI Tnf-Tnfrsf21 statistics (per Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.349936666666667
$ws.Range("H2").Value = 16.04981
$ws.Range("I2").Value = 0.04148245374559899
$ws.Range("J2").Value = 0.04148245374559899
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.465689
$ws.Range("N2").Value = 34.397067
$ws.Range("O2").Value = 0.1125836808441207
$ws.Range("P2").Value = 0.1279391038575984
$ws.Range("Q2").Value = 61.34070998969667
$ws.Range("R2").Value = 552.06638990727
$ws.Range("S2").Value = 0.004670247333125517
$ws.Range("T2").Value = 0.005307227958026211

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.349936666666667
$ws.Range("H3").Value = 16.04981
$ws.Range("I3").Value = 0.04148245374559899
$ws.Range("J3").Value = 0.04148245374559899
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.399706666666667
$ws.Range("N3").Value = 4.199120000000001
$ws.Range("O3").Value = 0.01374397374945266
$ws.Range("P3").Value = 0.01561853078317749
$ws.Range("Q3").Value = 7.48834201857778
$ws.Range("R3").Value = 67.39507816720001
$ws.Range("S3").Value = 0.0005701337553423967
$ws.Range("T3").Value = 0.0006478949807873743

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.349936666666667
$ws.Range("H4").Value = 16.04981
$ws.Range("I4").Value = 0.04148245374559899
$ws.Range("J4").Value = 0.04148245374559899
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.364677
$ws.Range("N4").Value = 79.094031
$ws.Range("O4").Value = 0.2588795475724425
$ws.Range("P4").Value = 0.2941884389917637
$ws.Range("Q4").Value = 141.0493521871234
$ws.Range("R4").Value = 1269.44416968411
$ws.Range("S4").Value = 0.01073895885785544
$ws.Range("T4").Value = 0.01220365831296581

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.349936666666667
$ws.Range("H5").Value = 16.04981
$ws.Range("I5").Value = 0.04148245374559899
$ws.Range("J5").Value = 0.04148245374559899
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.94195566666667
$ws.Range("N5").Value = 77.825867
$ws.Range("O5").Value = 0.2547287700938277
$ws.Range("P5").Value = 0.2894715320036049
$ws.Range("Q5").Value = 138.7878198261411
$ws.Range("R5").Value = 1249.09037843527
$ws.Range("S5").Value = 0.01056677442309053
$ws.Range("T5").Value = 0.01200798943700722

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.349936666666667
$ws.Range("H6").Value = 16.04981
$ws.Range("I6").Value = 0.04148245374559899
$ws.Range("J6").Value = 0.04148245374559899
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 36.6694545
$ws.Range("N6").Value = 73.338909
$ws.Range("O6").Value = 0.3600640277401564
$ws.Range("P6").Value = 0.2727823943638554
$ws.Range("Q6").Value = 196.179259176215
$ws.Range("R6").Value = 1177.07555505729
$ws.Range("S6").Value = 0.01493633937618511
$ws.Range("T6").Value = 0.01131568305681237

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.36117
$ws.Range("H7").Value = 4.08351
$ws.Range("I7").Value = 0.01055426915924182
$ws.Range("J7").Value = 0.01055426915924182
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 11.465689
$ws.Range("N7").Value = 34.397067
$ws.Range("O7").Value = 0.1125836808441207
$ws.Range("P7").Value = 0.1279391038575984
$ws.Range("Q7").Value = 15.60675189613
$ws.Range("R7").Value = 140.46076706517
$ws.Range("S7").Value = 0.001188238470567027
$ws.Range("T7").Value = 0.001350303738105287

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.36117
$ws.Range("H8").Value = 4.08351
$ws.Range("I8").Value = 0.01055426915924182
$ws.Range("J8").Value = 0.01055426915924182
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.399706666666667
$ws.Range("N8").Value = 4.199120000000001
$ws.Range("O8").Value = 0.01374397374945266
$ws.Range("P8").Value = 0.01561853078317749
$ws.Range("Q8").Value = 1.905238723466667
$ws.Range("R8").Value = 17.1471485112
$ws.Range("S8").Value = 0.0001450575982692774
$ws.Range("T8").Value = 0.0001648421777575592

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.36117
$ws.Range("H9").Value = 4.08351
$ws.Range("I9").Value = 0.01055426915924182
$ws.Range("J9").Value = 0.01055426915924182
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 26.364677
$ws.Range("N9").Value = 79.094031
$ws.Range("O9").Value = 0.2588795475724425
$ws.Range("P9").Value = 0.2941884389917637
$ws.Range("Q9").Value = 35.88680739209001
$ws.Range("R9").Value = 322.98126652881
$ws.Range("S9").Value = 0.002732284424902305
$ws.Range("T9").Value = 0.003104943968656265

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.36117
$ws.Range("H10").Value = 4.08351
$ws.Range("I10").Value = 0.01055426915924182
$ws.Range("J10").Value = 0.01055426915924182
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 25.94195566666667
$ws.Range("N10").Value = 77.825867
$ws.Range("O10").Value = 0.2547287700938277
$ws.Range("P10").Value = 0.2894715320036049
$ws.Range("Q10").Value = 35.31141179479668
$ws.Range("R10").Value = 317.80270615317
$ws.Range("S10").Value = 0.002688476002172886
$ws.Range("T10").Value = 0.003055160462704128

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.36117
$ws.Range("H11").Value = 4.08351
$ws.Range("I11").Value = 0.01055426915924182
$ws.Range("J11").Value = 0.01055426915924182
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 36.6694545
$ws.Range("N11").Value = 73.338909
$ws.Range("O11").Value = 0.3600640277401564
$ws.Range("P11").Value = 0.2727823943638554
$ws.Range("Q11").Value = 49.91336138176501
$ws.Range("R11").Value = 299.48016829059
$ws.Range("S11").Value = 0.003800212663330325
$ws.Range("T11").Value = 0.002879018812018579

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 44.98903266666667
$ws.Range("H12").Value = 134.967098
$ws.Range("I12").Value = 0.3488369270391816
$ws.Range("J12").Value = 0.3488369270391816
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 11.465689
$ws.Range("N12").Value = 34.397067
$ws.Range("O12").Value = 0.1125836808441207
$ws.Range("P12").Value = 0.1279391038575984
$ws.Range("Q12").Value = 515.8302569668407
$ws.Range("R12").Value = 4642.472312701567
$ws.Range("S12").Value = 0.03927334526042305
$ws.Range("T12").Value = 0.04462988383783133

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 44.98903266666667
$ws.Range("H13").Value = 134.967098
$ws.Range("I13").Value = 0.3488369270391816
$ws.Range("J13").Value = 0.3488369270391816
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.399706666666667
$ws.Range("N13").Value = 4.199120000000001
$ws.Range("O13").Value = 0.01374397374945266
$ws.Range("P13").Value = 0.01561853078317749
$ws.Range("Q13").Value = 62.97144895041779
$ws.Range("R13").Value = 566.7430405537601
$ws.Range("S13").Value = 0.004794405568066245
$ws.Range("T13").Value = 0.005448320283270498

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 44.98903266666667
$ws.Range("H14").Value = 134.967098
$ws.Range("I14").Value = 0.3488369270391816
$ws.Range("J14").Value = 0.3488369270391816
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 26.364677
$ws.Range("N14").Value = 79.094031
$ws.Range("O14").Value = 0.2588795475724425
$ws.Range("P14").Value = 0.2941884389917637
$ws.Range("Q14").Value = 1186.121314799115
$ws.Range("R14").Value = 10675.09183319204
$ws.Range("S14").Value = 0.09030674584846446
$ws.Range("T14").Value = 0.1026237910283406

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 44.98903266666667
$ws.Range("H15").Value = 134.967098
$ws.Range("I15").Value = 0.3488369270391816
$ws.Range("J15").Value = 0.3488369270391816
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 25.94195566666667
$ws.Range("N15").Value = 77.825867
$ws.Range("O15").Value = 0.2547287700938277
$ws.Range("P15").Value = 0.2894715320036049
$ws.Range("Q15").Value = 1167.103490924885
$ws.Range("R15").Value = 10503.93141832397
$ws.Range("S15").Value = 0.08885880138800105
$ws.Range("T15").Value = 0.1009783596894616

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 44.98903266666667
$ws.Range("H16").Value = 134.967098
$ws.Range("I16").Value = 0.3488369270391816
$ws.Range("J16").Value = 0.3488369270391816
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 36.6694545
$ws.Range("N16").Value = 73.338909
$ws.Range("O16").Value = 0.3600640277401564
$ws.Range("P16").Value = 0.2727823943638554
$ws.Range("Q16").Value = 1649.723286369347
$ws.Range("R16").Value = 9898.339718216084
$ws.Range("S16").Value = 0.1256036289742268
$ws.Range("T16").Value = 0.09515657220027748

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 77.26852533333333
$ws.Range("H17").Value = 231.805576
$ws.Range("I17").Value = 0.5991263500559777
$ws.Range("J17").Value = 0.5991263500559777
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 11.465689
$ws.Range("N17").Value = 34.397067
$ws.Range("O17").Value = 0.1125836808441207
$ws.Range("P17").Value = 0.1279391038575984
$ws.Range("Q17").Value = 885.9368809606212
$ws.Range("R17").Value = 7973.431928645592
$ws.Range("S17").Value = 0.06745184978000514
$ws.Range("T17").Value = 0.07665168832363559

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 77.26852533333333
$ws.Range("H18").Value = 231.805576
$ws.Range("I18").Value = 0.5991263500559777
$ws.Range("J18").Value = 0.5991263500559777
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 1.399706666666667
$ws.Range("N18").Value = 4.199120000000001
$ws.Range("O18").Value = 0.01374397374945266
$ws.Range("P18").Value = 0.01561853078317749
$ws.Range("Q18").Value = 108.1532700325689
$ws.Range("R18").Value = 973.3794302931202
$ws.Range("S18").Value = 0.008234376827774743
$ws.Range("T18").Value = 0.009357473341362062

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 77.26852533333333
$ws.Range("H19").Value = 231.805576
$ws.Range("I19").Value = 0.5991263500559777
$ws.Range("J19").Value = 0.5991263500559777
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 26.364677
$ws.Range("N19").Value = 79.094031
$ws.Range("O19").Value = 0.2588795475724425
$ws.Range("P19").Value = 0.2941884389917637
$ws.Range("Q19").Value = 2037.159712679651
$ws.Range("R19").Value = 18334.43741411686
$ws.Range("S19").Value = 0.1551015584412203
$ws.Range("T19").Value = 0.176256045681801

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 77.26852533333333
$ws.Range("H20").Value = 231.805576
$ws.Range("I20").Value = 0.5991263500559777
$ws.Range("J20").Value = 0.5991263500559777
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 25.94195566666667
$ws.Range("N20").Value = 77.825867
$ws.Range("O20").Value = 0.2547287700938277
$ws.Range("P20").Value = 0.2894715320036049
$ws.Range("Q20").Value = 2004.496658626044
$ws.Range("R20").Value = 18040.46992763439
$ws.Range("S20").Value = 0.1526147182805633
$ws.Range("T20").Value = 0.1734300224144319

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 77.26852533333333
$ws.Range("H21").Value = 231.805576
$ws.Range("I21").Value = 0.5991263500559777
$ws.Range("J21").Value = 0.5991263500559777
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 36.6694545
$ws.Range("N21").Value = 73.338909
$ws.Range("O21").Value = 0.3600640277401564
$ws.Range("P21").Value = 0.2727823943638554
$ws.Range("Q21").Value = 2833.394673992764
$ws.Range("R21").Value = 17000.36804395658
$ws.Range("S21").Value = 0.2157238467264142
$ws.Range("T21").Value = 0.163431120294747

